# Timesheet update: new clock-in/out entries for 2026-01-30, shift the
# running "Total Duration" summary down one row, and refresh row heights
# / selection to match the latest edit session.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: finish filling in the (previously half-empty) entry, but
#     change its date to 2026-01-30. Dates are stored as plain text in
#     this sheet, so force text formatting before typing the date-shaped
#     string to stop Excel from auto-converting it to a date serial.
$ws.Cells.Item(6, 1).NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "2026-01-30"
$ws.Cells.Item(6, 2).Value = "14:58:14"
$ws.Cells.Item(6, 3).Value = "15:13:42"
$ws.Cells.Item(6, 4).Value = "0.26 Hours"

# --- Row 7: a new clock-in/out pair for the same day, followed by the
#     "Total Duration" label/value that used to live on row 8.
$ws.Cells.Item(7, 1).NumberFormat = "@"
$ws.Cells.Item(7, 1).Value = "2026-01-30"
$ws.Cells.Item(7, 2).Value = "16:23:14"
$ws.Cells.Item(7, 3).Value = "Total Duration:"
$ws.Cells.Item(7, 4).Value = "11 Hours"

# Re-apply the normal (non-text) cell formatting used throughout the
# sheet to the two date cells we just typed into, so they pick up the
# same style as their neighbours instead of a one-off "text number
# format" style.
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(6, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(5, 1).Copy()
$ws.Cells.Item(7, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Row 8: the "Total Duration" summary moved up to row 7, so clear
#     the old C8/D8 text and bring A8/B8 into the formatted range (same
#     look as the rest of the table, just empty) to match.
$ws.Cells.Item(7, 3).Copy()
$ws.Cells.Item(8, 1).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(7, 3).Copy()
$ws.Cells.Item(8, 2).PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Cells.Item(8, 3).ClearContents()
$ws.Cells.Item(8, 4).ClearContents()
$excel.CutCopyMode = $false

# --- Row heights: small refresh pass across the body rows plus the two
#     newly-active rows at the bottom.
$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 15.65
$ws.Rows.Item(3).RowHeight = 15.65
$ws.Rows.Item(5).RowHeight = 15.65
$ws.Rows.Item(6).RowHeight = 15.65
$ws.Rows.Item(7).RowHeight = 15.5
$ws.Rows.Item(8).RowHeight = 15.5

# --- Selection follows the last entry the user touched.
$ws.Range("A7").Select() | Out-Null
